$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells I1 and J1 - match the style used by the other header cells (H1 etc.)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data rows 2-33 for columns I (I0) and J (IF)
$data = @{
    2  = @(9, 9)
    3  = @(9, 9)
    4  = @(6, 6)
    5  = @(6, 6)
    6  = @(6, 7)
    7  = @(7, 7)
    8  = @(5, 5)
    9  = @(9, 9)
    10 = @(4, 5)
    11 = @(4, 6)
    12 = @(8, 8)
    13 = @(8, 8)
    14 = @(8, 8)
    15 = @(7, 8)
    16 = @(6, 6)
    17 = @(6, 7)
    18 = @(8, 8)
    19 = @(4, 4)
    20 = @(7, 7)
    21 = @(7, 7)
    22 = @(9, 9)
    23 = @(8, 8)
    24 = @(7, 7)
    25 = @(6, 7)
    26 = @(7, 7)
    27 = @(6, 7)
    28 = @(9, 9)
    29 = @(2, 2)
    30 = @(5, 5)
    31 = @(3, 3)
    32 = @(4, 4)
    33 = @(9, 9)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
